# Auto-generated script: updates market price / profit columns (H:N)
# on each class sheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to reflect
# refreshed marketboard data pulled by the scheduled runner.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1030.8928
$ws.Range("I112").Value = 639.1667
$ws.Range("J112").Value = 1137.7273
$ws.Range("K112").Value = 1917.5001
$ws.Range("L112").Value = 3413.1819
$ws.Range("M112").Value = -809.5001
$ws.Range("N112").Value = -5629.1819
$ws.Range("H125").Value = 1680.6666
$ws.Range("I125").Value = 2009.1428
$ws.Range("J125").Value = 1545.4117
$ws.Range("K125").Value = 18082.2852
$ws.Range("L125").Value = 13908.7053
$ws.Range("M125").Value = -15622.2852
$ws.Range("N125").Value = -18828.7053
$ws.Range("H129").Value = 1130.7
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 1130.7
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 3392.1
$ws.Range("M129").ClearContents()
$ws.Range("N129").Value = -13392.1
$ws.Range("H138").Value = 1591.2307
$ws.Range("I138").Value = 1262.5758
$ws.Range("J138").Value = 3398.8333
$ws.Range("K138").Value = 3787.7274
$ws.Range("L138").Value = 10196.4999
$ws.Range("M138").Value = 1352.2726
$ws.Range("N138").Value = -20476.4999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19098.475
$ws.Range("I32").Value = 3589.311
$ws.Range("J32").Value = 218502
$ws.Range("K32").Value = 3589.311
$ws.Range("L32").Value = 218502
$ws.Range("M32").Value = -3302.311
$ws.Range("N32").Value = -219076
$ws.Range("H74").Value = 701.5333000000001
$ws.Range("I74").Value = 572.5
$ws.Range("J74").Value = 849
$ws.Range("K74").Value = 572.5
$ws.Range("L74").Value = 849
$ws.Range("M74").Value = 301.5
$ws.Range("N74").Value = -2597
$ws.Range("H77").Value = 701.5333000000001
$ws.Range("I77").Value = 572.5
$ws.Range("J77").Value = 849
$ws.Range("K77").Value = 2862.5
$ws.Range("L77").Value = 4245
$ws.Range("M77").Value = 1505.5
$ws.Range("N77").Value = -12981
$ws.Range("H97").Value = 1405.7142
$ws.Range("I97").Value = 964.93335
$ws.Range("J97").Value = 2507.6667
$ws.Range("K97").Value = 964.93335
$ws.Range("L97").Value = 2507.6667
$ws.Range("M97").Value = -468.93335
$ws.Range("N97").Value = -3499.6667
$ws.Range("H132").Value = 1504.5
$ws.Range("I132").Value = 1504.5403
$ws.Range("J132").Value = 1504
$ws.Range("K132").Value = 4513.6209
$ws.Range("L132").Value = 4512
$ws.Range("M132").Value = -1983.6209
$ws.Range("N132").Value = -9572
$ws.Range("H134").Value = 30000
$ws.Range("J134").Value = 30000
$ws.Range("L134").Value = 30000
$ws.Range("N134").Value = -40140

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H96").Value = 23000
$ws.Range("I96").Value = 13000
$ws.Range("J96").Value = 28000
$ws.Range("K96").Value = 13000
$ws.Range("L96").Value = 28000
$ws.Range("M96").Value = -10254
$ws.Range("N96").Value = -33492

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 388
$ws.Range("I22").Value = 388
$ws.Range("K22").Value = 388
$ws.Range("M22").Value = -38
$ws.Range("H31").Value = 30536.24
$ws.Range("I31").Value = 1458.4231
$ws.Range("J31").Value = 62037.207
$ws.Range("K31").Value = 1458.4231
$ws.Range("L31").Value = 62037.207
$ws.Range("M31").Value = -1163.4231
$ws.Range("N31").Value = -62627.207
$ws.Range("H34").Value = 30536.24
$ws.Range("I34").Value = 1458.4231
$ws.Range("J34").Value = 62037.207
$ws.Range("K34").Value = 1458.4231
$ws.Range("L34").Value = 62037.207
$ws.Range("M34").Value = -1256.4231
$ws.Range("N34").Value = -62441.207
$ws.Range("H99").Value = 13879
$ws.Range("I99").Value = 3566.8333
$ws.Range("J99").Value = 34503.332
$ws.Range("K99").Value = 3566.8333
$ws.Range("L99").Value = 34503.332
$ws.Range("M99").Value = -2068.8333
$ws.Range("N99").Value = -37499.332
$ws.Range("H122").Value = 578.8
$ws.Range("I122").Value = 545
$ws.Range("J122").Value = 714
$ws.Range("K122").Value = 1635
$ws.Range("L122").Value = 2142
$ws.Range("M122").Value = 815
$ws.Range("N122").Value = -7042
$ws.Range("H126").Value = 13879
$ws.Range("I126").Value = 3566.8333
$ws.Range("J126").Value = 34503.332
$ws.Range("K126").Value = 10700.4999
$ws.Range("L126").Value = 103509.996
$ws.Range("M126").Value = -8230.499899999999
$ws.Range("N126").Value = -108449.996
$ws.Range("H132").Value = 2111.6812
$ws.Range("I132").Value = 2000.537
$ws.Range("K132").Value = 6001.611
$ws.Range("M132").Value = -3471.611

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 603.4
$ws.Range("I34").Value = 393.75
$ws.Range("J34").Value = 743.1667
$ws.Range("K34").Value = 1181.25
$ws.Range("L34").Value = 2229.5001
$ws.Range("M34").Value = -1097.25
$ws.Range("N34").Value = -2397.5001
$ws.Range("H75").Value = 2305.5
$ws.Range("I75").Value = 622
$ws.Range("J75").Value = 2866.6667
$ws.Range("K75").Value = 1866
$ws.Range("L75").Value = 8600.000100000001
$ws.Range("M75").Value = -868
$ws.Range("N75").Value = -10596.0001
$ws.Range("H78").Value = 2305.5
$ws.Range("I78").Value = 622
$ws.Range("J78").Value = 2866.6667
$ws.Range("K78").Value = 5598
$ws.Range("L78").Value = 25800.0003
$ws.Range("M78").Value = -606
$ws.Range("N78").Value = -35784.0003
$ws.Range("H87").Value = 13000
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").ClearContents()
$ws.Range("H90").Value = 13000
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").ClearContents()
$ws.Range("H131").Value = 773.4299999999999
$ws.Range("I131").Value = 354.14285
$ws.Range("J131").Value = 841.68604
$ws.Range("K131").Value = 1062.42855
$ws.Range("L131").Value = 2525.05812
$ws.Range("M131").Value = 3977.57145
$ws.Range("N131").Value = -12605.05812
$ws.Range("H140").Value = 1293.4117
$ws.Range("I140").Value = 914.8261
$ws.Range("K140").Value = 2744.4783
$ws.Range("M140").Value = 2435.5217

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 61595.23
$ws.Range("I70").Value = 94969.77
$ws.Range("J70").Value = 5115.231
$ws.Range("K70").Value = 94969.77
$ws.Range("L70").Value = 5115.231
$ws.Range("M70").Value = -94699.77
$ws.Range("N70").Value = -5655.231
$ws.Range("H73").Value = 61595.23
$ws.Range("I73").Value = 94969.77
$ws.Range("J73").Value = 5115.231
$ws.Range("K73").Value = 94969.77
$ws.Range("L73").Value = 5115.231
$ws.Range("M73").Value = -94033.77
$ws.Range("N73").Value = -6987.231
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()
$ws.Range("H97").Value = 142859710
$ws.Range("I97").Value = 142859710
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 142859710
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -142859214
$ws.Range("N97").ClearContents()
$ws.Range("H102").Value = 233054.61
$ws.Range("I102").Value = 1774.3636
$ws.Range("J102").Value = 402660.12
$ws.Range("K102").Value = 1774.3636
$ws.Range("L102").Value = 402660.12
$ws.Range("M102").Value = -152.3635999999999
$ws.Range("N102").Value = -405904.12
$ws.Range("H122").Value = 501.85715
$ws.Range("I122").Value = 485.5
$ws.Range("J122").Value = 600
$ws.Range("K122").Value = 1456.5
$ws.Range("L122").Value = 1800
$ws.Range("M122").Value = 993.5
$ws.Range("N122").Value = -6700
$ws.Range("H126").Value = 8406586
$ws.Range("I126").Value = 3518.75
$ws.Range("J126").Value = 19610676
$ws.Range("K126").Value = 10556.25
$ws.Range("L126").Value = 58832028
$ws.Range("M126").Value = -8086.25
$ws.Range("N126").Value = -58836968
$ws.Range("H132").Value = 1906.758
$ws.Range("I132").Value = 1922.1818
$ws.Range("J132").Value = 1785.5714
$ws.Range("K132").Value = 5766.5454
$ws.Range("L132").Value = 5356.7142
$ws.Range("M132").Value = -3236.5454
$ws.Range("N132").Value = -10416.7142

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4187.636
$ws.Range("I7").Value = 1898.5
$ws.Range("J7").Value = 5495.7144
$ws.Range("K7").Value = 1898.5
$ws.Range("L7").Value = 5495.7144
$ws.Range("M7").Value = -1786.5
$ws.Range("N7").Value = -5719.7144
$ws.Range("H45").Value = 6181.381
$ws.Range("I45").Value = 4360.5
$ws.Range("K45").Value = 4360.5
$ws.Range("M45").Value = -3953.5
$ws.Range("H100").Value = 2425.4443
$ws.Range("I100").Value = 1980
$ws.Range("J100").Value = 2982.25
$ws.Range("K100").Value = 1980
$ws.Range("L100").Value = 2982.25
$ws.Range("M100").Value = -1439
$ws.Range("N100").Value = -4064.25
$ws.Range("H126").Value = 4187.636
$ws.Range("I126").Value = 1898.5
$ws.Range("J126").Value = 5495.7144
$ws.Range("K126").Value = 5695.5
$ws.Range("L126").Value = 16487.1432
$ws.Range("M126").Value = -3225.5
$ws.Range("N126").Value = -21427.1432

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 201534.1
$ws.Range("J81").Value = 144648.72
$ws.Range("L81").Value = 289297.44
$ws.Range("N81").Value = -291419.44
$ws.Range("H84").Value = 201534.1
$ws.Range("J84").Value = 144648.72
$ws.Range("L84").Value = 1446487.2
$ws.Range("N84").Value = -1457095.2
$ws.Range("H122").Value = 2450.889
$ws.Range("I122").Value = 2112.889
$ws.Range("J122").Value = 2788.889
$ws.Range("K122").Value = 6338.667
$ws.Range("L122").Value = 8366.667000000001
$ws.Range("M122").Value = -3888.667
$ws.Range("N122").Value = -13266.667
$ws.Range("H132").Value = 1766.1968
$ws.Range("I132").Value = 1814.6
$ws.Range("J132").Value = 1322.5
$ws.Range("K132").Value = 5443.799999999999
$ws.Range("L132").Value = 3967.5
$ws.Range("M132").Value = -2913.799999999999
$ws.Range("N132").Value = -9027.5
